# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets.
# These values mirror the same set of events, with the "全部类型" sheet
# offset by one extra row (it has one additional leading entry).

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll     = $wb.Worksheets.Item("全部类型")

# Row -> new F value for "展览" sheet
$exhibitUpdates = @{
    2  = 3108
    3  = 513
    5  = 71
    9  = 1094
    10 = 15267
    12 = 155
    14 = 6042
    17 = 61
    19 = 1256
    24 = 843
    27 = 128
    28 = 10903
    32 = 146
    33 = 3778
}

foreach ($row in $exhibitUpdates.Keys) {
    $wsExhibit.Range("F$row").Value = $exhibitUpdates[$row]
}

# Row -> new F value for "全部类型" sheet
$allUpdates = @{
    3  = 3108
    4  = 513
    6  = 71
    10 = 1094
    11 = 15267
    13 = 155
    15 = 6042
    18 = 61
    20 = 1256
    25 = 843
    28 = 128
    30 = 10903
    34 = 146
    35 = 3778
}

foreach ($row in $allUpdates.Keys) {
    $wsAll.Range("F$row").Value = $allUpdates[$row]
}
